# Add a new test case row (row 12) to the batch SQL test sheet, describing
# a case about auto_increment refresh and json \r\n escape (batch_011).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12 values --------------------------------------------------------
$ws.Range("A12").Value = "batch_011"
$ws.Range("B12").Value = "n"
$ws.Range("C12").Value = "批量操作语句11执行"
$ws.Range("D12").Value = "batchsql"
$ws.Range("E12").Value = "SingleTable"
$ws.Range("G12").Value = "batch011"
$ws.Range("I12").Value = "batch_sql_11"
$ws.Range("J12").Value = "select * from `$batch011"
$ws.Range("K12").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/batchsql/expectedresult/batch_011.csv"
$ws.Range("N12").Value = "csv_containsAll"

# --- Formatting: keep these cells as text (@) like the rest of the table --
$ws.Range("A12").NumberFormat = "@"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("I12").NumberFormat = "@"
$ws.Range("J12").NumberFormat = "@"
$ws.Range("K12").NumberFormat = "@"
$ws.Range("N12").NumberFormat = "@"

# Column K uses a "fill" horizontal alignment throughout the table.
$ws.Range("K12").HorizontalAlignment = 5

# Leave the selection where the author left it when saving the workbook.
$ws.Range("G22").Select() | Out-Null
